# 8-Bugtracker.xlsx edit
#
# The author inserted a new "Severidad" column between the existing
# "Impacto" and "Prioridad" columns (shifting "Prioridad" and
# "Fecha solución" one column to the right), and renamed the
# "Descripción" header (column C) to "Reporte" (the underlying long-form
# description text in C2:C4 is unchanged).
#
# Concretely, on the INCIDENTES sheet:
#   - Column C header: "Descripción" -> "Reporte"
#   - A new blank column is inserted at K (header "Severidad"); the old
#     K column ("Prioridad") becomes L, and the old L column
#     ("Fecha solución") becomes M.
#   - Selection ends up on K2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the old "Prioridad" column (K),
# pushing Prioridad -> L and Fecha solución -> M.
$ws.Columns("K:K").Insert() | Out-Null

# Rename the "Descripción" header to "Reporte".
$ws.Range("C1").Value = "Reporte"

# Header for the newly inserted column.
$ws.Range("K1").Value = "Severidad"

# Narrow column B slightly (author re-fit it to its contents after the
# layout change); closest value this runtime's width-rounding allows.
$ws.Columns("B:B").ColumnWidth = 26.8

# Match the author's final selection.
$ws.Range("K2").Select() | Out-Null
